# Applies the OOXML change described by the diff:
#  1. Slide 2, shape "CustomShape 7" (the "Trabalhadores envolvidos" bullet
#     list box): remove the bullet paragraph
#     "Verifica a disponibilidade do produtos." entirely. The shape has
#     <a:spAutoFit/>, so removing the paragraph's text (run + trailing
#     paragraph mark) also shrinks the shape's height to match the
#     target extent (cy 1891372 -> 1752872 EMU).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(7)

$tr = $sh.TextFrame.TextRange
$fullText = $tr.Text

$needle = "Verifica a disponibilidade do produtos."
$idx = $fullText.IndexOf($needle)

if ($idx -ge 0) {
    $delLen = $needle.Length

    # Pull in the trailing paragraph-mark (carriage return) too, so the
    # whole paragraph -- not just its text -- disappears, and subsequent
    # paragraphs keep their own original formatting.
    $endPos = $idx + $delLen
    if ($endPos -lt $fullText.Length -and $fullText.Substring($endPos, 1) -eq "`r") {
        $delLen = $delLen + 1
    }

    $victim = $tr.Characters($idx + 1, $delLen)
    $victim.Delete()
}
